$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Crime complaints table updates (rows 14-30) ---
$updates = @{
    14 = @{ C=4;   F=7;    G=1;   H=600;                   I=43;   K=34.375;              L=-4.444444444444;    M=26.470588235294;    N=-81.057268722467 }
    15 = @{ C=3;   D=4;    E=-25; F=13;   G=15;  H=-13.333333333333; I=91;   J=133;  K=-31.578947368421;   L=-7.142857142857;    M=-34.532374100719;   N=-73.699421965317 }
    16 = @{ C=44;  D=31;   E=41.935483870967; F=168;  G=175; H=-4;               I=1294; J=1421; K=-8.937368050668;    L=12.326388888888;    M=-17.945466074825;   N=-79.953524399690 }
    17 = @{ C=57;  D=49;   E=16.326530612244; F=235;  G=237; H=-0.843881856540;  I=2036; J=2104; K=-3.231939163498;    L=9.227467811158;     M=52.281226626776;    N=-49.753208292201 }
    18 = @{ C=28;  D=35;   E=-20; F=103;  G=153; H=-32.679738562091; I=1022; J=1234; K=-17.179902755267;   L=10.966340933767;    M=9.892473118279;     N=-86.331416343453 }
    19 = @{ C=137; D=136;  E=0.735294117647;  F=587;  G=522; H=12.452107279693; I=4292; J=4368; K=-1.739926739926;    L=29.628511023859;    M=37.564102564102;    N=-41.952934812009 }
    20 = @{ C=18;  D=28;   E=-35.714285714285; F=100; G=84;  H=19.047619047619; I=913;  J=873;  K=4.581901489117;     L=43.103448275862;    M=135.917312661499;   N=-85.592551680606 }
    21 = @{ C=291; D=283;  E=2.826855123674;  F=1213; G=1187; H=2.190395956192;  I=9691; J=10165; K=-4.663059517953;   L=20.699962635446;    M=28.801169590643;    N=-69.985753221010 }
    22 = @{ C=5;   D=9;    E=-44.444444444444; F=24;  G=21;  H=14.285714285714; I=193;  J=211;  K=-8.530805687203;    L=25.324675324675;    M=23.717948717948 }
    23 = @{ C=22;  E=-4.347826086956; F=109; G=100; H=9;                I=859;  J=873;  K=-1.603665521191;    L=5.012224938875;     M=50.966608084358 }
    24 = @{ C=280; D=340;  E=-17.647058823529; F=1106; G=1343; H=-17.647058823529; I=9564; J=10742; K=-10.966300502699;  L=17.292126563649;    M=53.663239074550 }
    25 = @{ C=84;  D=90;   E=-6.666666666666; F=321;  G=333; H=-3.603603603603;  I=3062; J=3073; K=-0.357956394402;    L=10.821570756424;    M=-16.521264994547 }
    26 = @{ C=6;   D=9;    E=-33.333333333333; F=22;  G=23;  H=-4.347826086956; I=171;  J=207;  K=-17.391304347826;   L=-0.581395348837 }
    27 = @{ C=8;   D=13;   E=-38.461538461538; F=39;  H=-17.021276595744;        I=401;  J=461;  K=-13.015184381778;   L=-8.237986270022 }
    28 = @{ C=3;   D=1;    E=200; F=13;  G=19;  H=-31.578947368421; I=113;  J=140;  K=-19.285714285714;   L=-35.428571428571;   M=-19.858156028368;   N=-78.917910447761 }
    29 = @{ C=2;   D=1;    E=100; F=12;  G=17;  H=-29.411764705882; I=101;  J=117;  K=-13.675213675213;   L=-33.552631578947;   M=-18.548387096774;   N=-79.554655870445 }
    30 = @{ D=4;   G=9;    H=-44.444444444444; J=75; K=-36 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
